# add OnResume event processing
# Update prompt text for the ON_INTRO, RED_PICK_UP_IPAD, BLUE_PICK_UP_IPAD,
# READ_THE_WORD and CORRECT_SELECTION rows in the prompts table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = "This game will help you and your partner learn to read. One of you will be the red partner and one will be the blue partner. You decide together who will be blue and who will be red!"
$ws.Range("H5").Value = "Red partner, pick up the iPad and hold it in front of your face so your partner is looking at the back of the iPad."
$ws.Range("H6").Value = "Blue partner,  pick up the iPad and hold it in front of your face so your partner is looking at the back of the iPad."
$ws.Range("H7").Value = "OK, Now read the word you see out loud so your partner can hear it. When you have read the word, place the iPad flat in front of you again."
$ws.Range("H13").Value = "You're right! The answer was [ANSWER]"

$ws.Range("H13").Select()
